$d = $word.ActiveDocument

# --- Edit 1 ---
# "This book provides a framework for using decision making as the basis for
# understanding organizations.  It provides a useful framework for structuring
# a study of the role of development stage in university technology transfer."
# becomes
# "This source provides a useful scaffold for structuring a study of the role
# of development stage in university technology transfer.  It provides a
# framework for using decision making as the basis for understanding
# organizations.  Chapters 4, 6, 8, and 10 focus on the sociology of
# administration ... (pp. 356-360)."
$old1 = "This book provides a framework for using decision making as the basis for understanding organizations.  It provides a useful framework for structuring a study of the role of development stage in university technology transfer.  Simon argued that decision making"
$new1 = "This source provides a useful scaffold for structuring a study of the role of development stage in university technology transfer.  It provides a framework for using decision making as the basis for understanding organizations.  Chapters 4, 6, 8, and 10 focus on the sociology of administration " + [char]0x2013 + " what one might call descriptive administration theory.  Chapters 3, 9, and 11 emphasize what Simon calls the practical science of administration " + [char]0x2013 + " what might be apply labeled as normative administration theory (pp. 356-360).  Simon argued that decision making"

$result1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Edit 1 result:" $result1

# --- Edit 2 ---
# Remove the "Vertical decision making..." sentence and insert the new
# "As with physical tasks, there is specialization..." sentence plus "general".
$old2 = "the organization hierarchy (p. 2).  Vertical decision making refers to the division of decision making responsibilities between operative and supervisory personnel within the organization (p. 23).  Simon argued that two kinds of decisions are made in organizations (p. 4)."
$new2 = "the organization hierarchy (p. 2).  As with physical tasks, there is specialization regarding decisions in organizations.  Simon argued that two general kinds of decisions are made in organizations (p. 4)."

$result2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Host "Edit 2 result:" $result2

# --- Edit 3 ---
# Append the new role-theory commentary at the end of the paragraph.
$old3 = "the organization and are assembled into a decision."
$new3 = "the organization and are assembled into a final decision.  Also relevant is Simon" + [char]0x2019 + "s critique of role theory and the idea that roles determine behavior, which he argues is too constraining in its original connotation of dramatic part.  Simon counters that a role specifies some, but not all, of the premises that underlie a decision (pp. 24-25)."

$result3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Host "Edit 3 result:" $result3

# --- Edit 4 ---
# The added text reflows the document by one page, so the cached page-number
# field in the running header (used on all pages but the title page) drops
# from 5 to 4.
$sec = $d.Sections.Item(2)
$header = $sec.Headers.Item(1)
$result4 = $header.Range.Find.Execute("5", $true, $false, $false, $false, $false, $true, 1, $false, "4", 2)
Write-Host "Edit 4 result:" $result4

